# Daily auto-push edit: insert one new sampling row ("2026/01/06") into the
# sei1 (Sheet1) log just before the existing "2026/12/29" block, shifting
# every row from the old 565 down through 606 to 566..607, and growing the
# sheet's used range from D606 to D607.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift rows 565..606 down to 566..607 by inserting a blank row at 565.
$ws.Rows.Item(565).Insert()

# Column A holds plain text dates (e.g. "2026/12/29"), not real Excel date
# serials. Force the new cell to text first so the "YYYY/MM/DD" string isn't
# auto-parsed into a date value, then drop the format back to Normal so no
# stray number-formatting is left behind on the cell.
$ws.Range("A565").NumberFormat = "@"
$ws.Range("A565").Value = "2026/01/06"
$ws.Range("A565").Style = "Normal"

$ws.Range("B565").Value = "火"
$ws.Range("C565").Value = 3
$ws.Range("D565").Value = 153
